# Doing Updates for Financials
# The quarterly financials sheet gains two new (most-recent) quarter columns.
# Existing data in columns D:K (8 quarters) shifts right to F:M, and the two
# newly inserted columns D:E are populated with the latest two quarters of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank columns before column D; this shifts the existing
# D:K data to F:M automatically.
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy number/date formatting from the (now shifted) old column D/E — which
# live at F/G — into the freshly inserted, still-unformatted D:E columns so
# the new cells pick up the same date / numeric styles used throughout the
# report (row 7 onward holds data; rows 5-6 are plain text labels with no
# per-quarter cells).
$ws.Range("F7:G102").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New values for the two inserted quarter columns (D = newest quarter,
# E = second-newest quarter), keyed by worksheet row number. A $null entry
# means the row has no data in that quarter (blank, formatted cell only).
$newData = @{
    7 = @(43465, 43373)
    8 = @(9300, 9000)
    9 = @("NA", "NA")
    10 = @("NA", "NA")
    11 = @($null, $null)
    12 = @(15100, 15100)
    13 = @(0, 0)
    14 = @(0, 0)
    15 = @(0, 0)
    16 = @($null, $null)
    17 = @(20700, 20500)
    18 = @(-11400, -11500)
    19 = @($null, $null)
    20 = @(1100, 1000)
    21 = @(-10200, -10300)
    22 = @(400, 400)
    23 = @(-10800, -10900)
    24 = @(0, 0)
    25 = @(0, 0)
    26 = @(-10800, -10900)
    27 = @(-10800, -10900)
    28 = @(0, 0)
    29 = @(0, 0)
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-1100, -1000)
    33 = @(-10800, -10900)
    34 = @(0, 0)
    35 = @(-10800, -10900)
    38 = @(43465, 43373)
    39 = @($null, $null)
    40 = @($null, $null)
    41 = @(28100, 24100)
    42 = @(148900, 149000)
    43 = @(2100, 300)
    44 = @(0, 0)
    45 = @(2300, 2600)
    46 = @(181400, 176100)
    47 = @(0, 12900)
    48 = @(1500, 1600)
    49 = @(0, 0)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(400, 200)
    53 = @(0, 0)
    54 = @(183300, 190800)
    55 = @($null, $null)
    56 = @($null, $null)
    57 = @(1000, 600)
    58 = @(0, 0)
    59 = @(63400, 60800)
    60 = @(64400, 61300)
    61 = @(19700, 14700)
    62 = @(84500, 93200)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(168600, 169300)
    67 = @($null, $null)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(-374500, -363700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(14700, 21500)
    77 = @(0, 0)
    80 = @(43465, 43373)
    81 = @(-10800, -10900)
    82 = @($null, $null)
    83 = @(100, 100)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(-15300, -17500)
    90 = @($null, $null)
    91 = @(-100, -100)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(13400, -28500)
    95 = @($null, $null)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(5900, 1400)
    101 = @(0, 0)
    102 = @(4000, -44600)
}

foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $dVal = $vals[0]
    $eVal = $vals[1]
    if ($dVal -ne $null) {
        $ws.Cells.Item($r, 4).Value = $dVal
    }
    if ($eVal -ne $null) {
        $ws.Cells.Item($r, 5).Value = $eVal
    }
}
